# Updates per-row profit/price figures (columns H-N) across several
# worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), per scheduled
# runner recalculation.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1782.2128
$ws.Range("I15").Value = 1782.2128
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 5346.6384
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -5177.6384

$ws.Range("H17").Value = 2273175
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2273175
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6819525
$ws.Range("N17").Value = -6819861

$ws.Range("H94").Value = 2532.1875
$ws.Range("I94").Value = 2322.5
$ws.Range("J94").Value = 4000
$ws.Range("K94").Value = 2322.5
$ws.Range("L94").Value = 4000
$ws.Range("M94").Value = -1871.5
$ws.Range("N94").Value = -4902

$ws.Range("H99").Value = 490.92856
$ws.Range("I99").Value = 232.44444
$ws.Range("J99").Value = 956.2
$ws.Range("K99").Value = 697.33332
$ws.Range("L99").Value = 2868.6
$ws.Range("M99").Value = 800.66668
$ws.Range("N99").Value = -5864.6

$ws.Range("H125").Value = 1536.3334
$ws.Range("I125").Value = 750
$ws.Range("J125").Value = 1798.4445
$ws.Range("K125").Value = 6750
$ws.Range("L125").Value = 16186.0005
$ws.Range("M125").Value = -4290
$ws.Range("N125").Value = -21106.0005

$ws.Range("H137").Value = 5715462.5
$ws.Range("I137").Value = 1214.625
$ws.Range("J137").Value = 10527461
$ws.Range("K137").Value = 3643.875
$ws.Range("L137").Value = 31582383
$ws.Range("M137").Value = -1093.875
$ws.Range("N137").Value = -31587483

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3564.5232
$ws.Range("I32").Value = 3259.5806
$ws.Range("J32").Value = 9866.666999999999
$ws.Range("K32").Value = 3259.5806
$ws.Range("L32").Value = 9866.666999999999
$ws.Range("M32").Value = -2972.5806

$ws.Range("H45").Value = 833.3333
$ws.Range("I45").Value = 750
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 750
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -373
$ws.Range("N45").Value = -1754

$ws.Range("H61").Value = 994.03845
$ws.Range("I61").Value = 692.10254
$ws.Range("J61").Value = 1899.8462
$ws.Range("K61").Value = 692.10254
$ws.Range("L61").Value = 1899.8462
$ws.Range("M61").Value = -480.10254

$ws.Range("H74").Value = 672
$ws.Range("I74").Value = 506.58823
$ws.Range("J74").Value = 1375
$ws.Range("K74").Value = 506.58823
$ws.Range("L74").Value = 1375
$ws.Range("M74").Value = 367.41177

$ws.Range("H77").Value = 672
$ws.Range("I77").Value = 506.58823
$ws.Range("J77").Value = 1375
$ws.Range("K77").Value = 2532.94115
$ws.Range("L77").Value = 6875
$ws.Range("M77").Value = 1835.05885

$ws.Range("H94").Value = 30000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 30000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31802

$ws.Range("H122").Value = 16703.143
$ws.Range("I122").Value = 2782.4
$ws.Range("J122").Value = 51505
$ws.Range("K122").Value = 8347.200000000001
$ws.Range("L122").Value = 154515
$ws.Range("M122").Value = -5897.200000000001

$ws.Range("H123").Value = 49014.5
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 49014.5
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 49014.5
$ws.Range("N123").Value = -58814.5

$ws.Range("H136").Value = 994.03845
$ws.Range("I136").Value = 692.10254
$ws.Range("J136").Value = 1899.8462
$ws.Range("K136").Value = 2076.30762
$ws.Range("L136").Value = 5699.5386
$ws.Range("M136").Value = 473.69238

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 530.88
$ws.Range("I99").Value = 459.5238
$ws.Range("J99").Value = 905.5
$ws.Range("K99").Value = 459.5238
$ws.Range("L99").Value = 905.5
$ws.Range("M99").Value = 1038.4762
$ws.Range("N99").Value = -3901.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1363.7858
$ws.Range("I16").Value = 900.1667
$ws.Range("J16").Value = 1711.5
$ws.Range("K16").Value = 900.1667
$ws.Range("L16").Value = 1711.5
$ws.Range("M16").Value = -613.1667
$ws.Range("N16").Value = -2285.5

$ws.Range("H31").Value = 5407758.5
$ws.Range("I31").Value = 2626.9048
$ws.Range("J31").Value = 12501994
$ws.Range("K31").Value = 2626.9048
$ws.Range("L31").Value = 12501994
$ws.Range("M31").Value = -2331.9048
$ws.Range("N31").Value = -12502584

$ws.Range("H34").Value = 5407758.5
$ws.Range("I34").Value = 2626.9048
$ws.Range("J34").Value = 12501994
$ws.Range("K34").Value = 2626.9048
$ws.Range("L34").Value = 12501994
$ws.Range("M34").Value = -2424.9048
$ws.Range("N34").Value = -12502398

$ws.Range("H54").Value = 21975
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 21975
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 21975
$ws.Range("N54").Value = -23291

$ws.Range("H74").Value = 32098.334
$ws.Range("I74").Value = 25000
$ws.Range("J74").Value = 33518
$ws.Range("K74").Value = 25000
$ws.Range("L74").Value = 33518
$ws.Range("M74").Value = -24126
$ws.Range("N74").Value = -35266

$ws.Range("H77").Value = 32098.334
$ws.Range("I77").Value = 25000
$ws.Range("J77").Value = 33518
$ws.Range("K77").Value = 75000
$ws.Range("L77").Value = 100554
$ws.Range("M77").Value = -70632
$ws.Range("N77").Value = -109290

$ws.Range("H88").Value = 30890.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 30890.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 30890.5
$ws.Range("N88").Value = -31702.5

$ws.Range("H91").Value = 30890.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 30890.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 30890.5
$ws.Range("N91").Value = -33698.5

$ws.Range("H94").Value = 2179.913
$ws.Range("I94").Value = 1449
$ws.Range("J94").Value = 2569.7334
$ws.Range("K94").Value = 1449
$ws.Range("L94").Value = 2569.7334
$ws.Range("M94").Value = -998
$ws.Range("N94").Value = -3471.7334

$ws.Range("H99").Value = 2270.6956
$ws.Range("I99").Value = 1879.4286
$ws.Range("J99").Value = 2879.3333
$ws.Range("K99").Value = 1879.4286
$ws.Range("L99").Value = 2879.3333
$ws.Range("M99").Value = -381.4286
$ws.Range("N99").Value = -5875.3333

$ws.Range("H113").Value = 1363.7858
$ws.Range("I113").Value = 900.1667
$ws.Range("J113").Value = 1711.5
$ws.Range("K113").Value = 900.1667
$ws.Range("L113").Value = 1711.5
$ws.Range("M113").Value = 1269.8333
$ws.Range("N113").Value = -6051.5

$ws.Range("H126").Value = 2270.6956
$ws.Range("I126").Value = 1879.4286
$ws.Range("J126").Value = 2879.3333
$ws.Range("K126").Value = 5638.2858
$ws.Range("L126").Value = 8637.999899999999
$ws.Range("M126").Value = -3168.2858
$ws.Range("N126").Value = -13577.9999

$ws.Range("H132").Value = 2311.7354
$ws.Range("I132").Value = 1991.7084
$ws.Range("J132").Value = 3079.8
$ws.Range("K132").Value = 5975.1252
$ws.Range("L132").Value = 9239.400000000001
$ws.Range("M132").Value = -3445.1252
$ws.Range("N132").Value = -14299.4

$ws.Range("H141").Value = 52715.23
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 52715.23
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 52715.23
$ws.Range("N141").Value = -63075.23

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1367.2184
$ws.Range("I68").Value = 1333.6757
$ws.Range("J68").Value = 1392.04
$ws.Range("K68").Value = 4001.0271
$ws.Range("L68").Value = 4176.12
$ws.Range("M68").Value = -3190.0271
$ws.Range("N68").Value = -5798.12

$ws.Range("H71").Value = 1367.2184
$ws.Range("I71").Value = 1333.6757
$ws.Range("J71").Value = 1392.04
$ws.Range("K71").Value = 12003.0813
$ws.Range("L71").Value = 12528.36
$ws.Range("M71").Value = -7947.0813
$ws.Range("N71").Value = -20640.36

$ws.Range("H81").Value = 6672.909
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 6672.909
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 20018.727
$ws.Range("N81").Value = -22264.727
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 6672.909
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 6672.909
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 60056.181
$ws.Range("N84").Value = -71288.181
$ws.Range("M84").ClearContents()

$ws.Range("H131").Value = 3328986.5
$ws.Range("I131").Value = 5782.5
$ws.Range("J131").Value = 6350081
$ws.Range("K131").Value = 17347.5
$ws.Range("L131").Value = 19050243
$ws.Range("M131").Value = -12307.5
$ws.Range("N131").Value = -19060323

$ws.Range("H132").Value = 2250
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -6470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H43").Value = 6730.4443
$ws.Range("I43").Value = 3443.5
$ws.Range("J43").Value = 9360
$ws.Range("K43").Value = 3443.5
$ws.Range("L43").Value = 9360
$ws.Range("M43").Value = -3292.5
$ws.Range("N43").Value = -9662

$ws.Range("H110").Value = 60000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 60000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 60000
$ws.Range("N110").Value = -68180

$ws.Range("H122").Value = 6219.7334
$ws.Range("I122").Value = 11046.667
$ws.Range("J122").Value = 3001.7778
$ws.Range("K122").Value = 33140.001
$ws.Range("L122").Value = 9005.3334
$ws.Range("M122").Value = -30690.001
$ws.Range("N122").Value = -13905.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 614.4
$ws.Range("I16").Value = 614.4
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 614.4
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -444.4

$ws.Range("H122").Value = 2855.4443
$ws.Range("I122").Value = 2004
$ws.Range("J122").Value = 2961.875
$ws.Range("K122").Value = 6012
$ws.Range("L122").Value = 8885.625
$ws.Range("M122").Value = -3562
$ws.Range("N122").Value = -13785.625

$ws.Range("H136").Value = 5910.88
$ws.Range("I136").Value = 8762.286
$ws.Range("J136").Value = 2281.818
$ws.Range("K136").Value = 26286.858
$ws.Range("L136").Value = 6845.454000000001
$ws.Range("M136").Value = -23736.858
$ws.Range("N136").Value = -11945.454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1312.8
$ws.Range("I107").Value = 1718.5714
$ws.Range("J107").Value = 366
$ws.Range("K107").Value = 5155.7142
$ws.Range("L107").Value = 1098
$ws.Range("M107").Value = -3235.7142
$ws.Range("N107").Value = -4938

$ws.Range("H122").Value = 8500
$ws.Range("I122").Value = 8500
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 25500
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -23050
$ws.Range("N122").ClearContents()

$ws.Range("H123").Value = 36581.785
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 36581.785
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 36581.785
$ws.Range("N123").Value = -46381.785

$ws.Range("H132").Value = 1040.8422
$ws.Range("I132").Value = 876.4722
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 2629.4166
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -99.41660000000002
$ws.Range("N132").Value = -17058.5

$ws.Range("H136").Value = 1894.7667
$ws.Range("I136").Value = 2233.5789
$ws.Range("J136").Value = 1309.5454
$ws.Range("K136").Value = 6700.736699999999
$ws.Range("L136").Value = 3928.6362
$ws.Range("M136").Value = -4150.736699999999
$ws.Range("N136").Value = -9028.636200000001
